$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7877616882324219
$ws.Range("B1").Value = 1.334078669548035
$ws.Range("C1").Value = 4.284316062927246
$ws.Range("D1").Value = 1.584183812141418
$ws.Range("E1").Value = 0.7181151509284973
